$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 458729.22
$ws.Range("I19").Value = 718476.4
$ws.Range("K19").Value = 718476.4
$ws.Range("M19").Value = -718301.4
$ws.Range("H32").Value = 1831.25
$ws.Range("J32").Value = 1978.5714
$ws.Range("L32").Value = 1978.5714
$ws.Range("N32").Value = -2630.5714
$ws.Range("H128").Value = 40000
$ws.Range("J128").Value = 40000
$ws.Range("L128").Value = 40000
$ws.Range("N128").Value = -49960
$ws.Range("H131").Value = 3263.5
$ws.Range("I131").Value = 613.75
$ws.Range("J131").Value = 4020.5715
$ws.Range("K131").Value = 1841.25
$ws.Range("L131").Value = 12061.7145
$ws.Range("M131").Value = 3198.75
$ws.Range("N131").Value = -22141.7145
$ws.Range("H135").Value = 469.89474
$ws.Range("I135").Value = 376.8125
$ws.Range("J135").Value = 966.3333
$ws.Range("K135").Value = 3391.3125
$ws.Range("L135").Value = 8696.9997
$ws.Range("M135").Value = -856.3125
$ws.Range("N135").Value = -13766.9997
$ws.Range("H137").Value = 1757.174
$ws.Range("I137").Value = 1439.2142
$ws.Range("J137").Value = 2251.7778
$ws.Range("K137").Value = 4317.642599999999
$ws.Range("L137").Value = 6755.3334
$ws.Range("M137").Value = -1767.642599999999
$ws.Range("N137").Value = -11855.3334
$ws.Range("H138").Value = 3525.6216
$ws.Range("J138").Value = 2646.52
$ws.Range("L138").Value = 7939.559999999999
$ws.Range("N138").Value = -18219.56
$ws.Range("H141").Value = 877557
$ws.Range("I141").Value = 1219639.6
$ws.Range("K141").Value = 3658918.8
$ws.Range("M141").Value = -3653738.8

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 2644.5269
$ws.Range("I32").Value = 2133.6707
$ws.Range("J32").Value = 6452.727
$ws.Range("K32").Value = 2133.6707
$ws.Range("L32").Value = 6452.727
$ws.Range("M32").Value = -1846.6707
$ws.Range("N32").Value = -7026.727
$ws.Range("H61").Value = 2695.158
$ws.Range("I61").Value = 963.4545000000001
$ws.Range("K61").Value = 963.4545000000001
$ws.Range("M61").Value = -751.4545000000001
$ws.Range("H136").Value = 2695.158
$ws.Range("I136").Value = 963.4545000000001
$ws.Range("K136").Value = 2890.3635
$ws.Range("M136").Value = -340.3635000000004

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H103").Value = 38888
$ws.Range("J103").Value = 38888
$ws.Range("L103").Value = 38888
$ws.Range("N103").Value = -41232

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1485.6154
$ws.Range("I31").Value = 846.7778
$ws.Range("J31").Value = 1823.8235
$ws.Range("K31").Value = 846.7778
$ws.Range("L31").Value = 1823.8235
$ws.Range("M31").Value = -551.7778
$ws.Range("N31").Value = -2413.8235
$ws.Range("H34").Value = 1485.6154
$ws.Range("I34").Value = 846.7778
$ws.Range("J34").Value = 1823.8235
$ws.Range("K34").Value = 846.7778
$ws.Range("L34").Value = 1823.8235
$ws.Range("M34").Value = -644.7778
$ws.Range("N34").Value = -2227.8235

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 721.8182
$ws.Range("J2").Value = 694
$ws.Range("L2").Value = 4164
$ws.Range("N2").Value = -4390
$ws.Range("H37").Value = 84500
$ws.Range("J37").Value = 84500
$ws.Range("L37").Value = 253500
$ws.Range("N37").Value = -253724
$ws.Range("H50").Value = 71479400
$ws.Range("I50").Value = 100197.29
$ws.Range("J50").Value = 142858610
$ws.Range("K50").Value = 300591.87
$ws.Range("L50").Value = 428575830
$ws.Range("M50").Value = -300110.87
$ws.Range("N50").Value = -428576792
$ws.Range("H53").Value = 71479400
$ws.Range("I53").Value = 100197.29
$ws.Range("J53").Value = 142858610
$ws.Range("K53").Value = 300591.87
$ws.Range("L53").Value = 428575830
$ws.Range("M53").Value = -300110.87
$ws.Range("N53").Value = -428576792
$ws.Range("H61").Value = 217.25
$ws.Range("J61").Value = 227.6
$ws.Range("L61").Value = 682.8
$ws.Range("N61").Value = -1112.8
$ws.Range("H68").Value = 2297.205
$ws.Range("I68").Value = 1290.1818
$ws.Range("J68").Value = 2692.8215
$ws.Range("K68").Value = 3870.5454
$ws.Range("L68").Value = 8078.4645
$ws.Range("M68").Value = -3059.5454
$ws.Range("N68").Value = -9700.4645
$ws.Range("H71").Value = 2297.205
$ws.Range("I71").Value = 1290.1818
$ws.Range("J71").Value = 2692.8215
$ws.Range("K71").Value = 11611.6362
$ws.Range("L71").Value = 24235.3935
$ws.Range("M71").Value = -7555.636200000001
$ws.Range("N71").Value = -32347.3935
$ws.Range("H107").Value = 1510.125
$ws.Range("I107").Value = 1750
$ws.Range("J107").Value = 1483.4722
$ws.Range("K107").Value = 5250
$ws.Range("L107").Value = 4450.4166
$ws.Range("M107").Value = -3330
$ws.Range("N107").Value = -8290.4166
$ws.Range("H122").Value = 1174
$ws.Range("J122").Value = 1437
$ws.Range("L122").Value = 12933
$ws.Range("N122").Value = -17833
$ws.Range("H131").Value = 11645295
$ws.Range("J131").Value = 18210.756
$ws.Range("L131").Value = 54632.268
$ws.Range("N131").Value = -64712.268
$ws.Range("H140").Value = 1282.6818
$ws.Range("I140").Value = 867.619
$ws.Range("K140").Value = 2602.857
$ws.Range("M140").Value = 2577.143

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 2609.9285
$ws.Range("I122").Value = 1671.2858
$ws.Range("J122").Value = 3548.5715
$ws.Range("K122").Value = 5013.857400000001
$ws.Range("L122").Value = 10645.7145
$ws.Range("M122").Value = -2563.857400000001
$ws.Range("N122").Value = -15545.7145

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H122").Value = 7288.706
$ws.Range("I122").Value = 8572.714
$ws.Range("J122").Value = 6389.9
$ws.Range("K122").Value = 25718.142
$ws.Range("L122").Value = 19169.7
$ws.Range("M122").Value = -23268.142
$ws.Range("N122").Value = -24069.7
$ws.Range("H132").Value = 4266.2856
$ws.Range("I132").Value = 1124.5
$ws.Range("K132").Value = 3373.5
$ws.Range("M132").Value = -843.5
$ws.Range("H136").Value = 3440.6316
$ws.Range("I136").Value = 2365.5386
$ws.Range("K136").Value = 7096.6158
$ws.Range("M136").Value = -4546.6158

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 157876.92
$ws.Range("I122").Value = 157876.92
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 473630.76
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = -471180.76
$ws.Range("N122").ClearContents()
$ws.Range("H123").Value = 38024.332
$ws.Range("J123").Value = 38024.332
$ws.Range("L123").Value = 38024.332
$ws.Range("N123").Value = -47824.332
$ws.Range("H126").Value = 5130.136
$ws.Range("I126").Value = 4950.9443
$ws.Range("K126").Value = 14852.8329
$ws.Range("M126").Value = -12382.8329
$ws.Range("H132").Value = 600.38464
$ws.Range("I132").Value = 492
$ws.Range("J132").Value = 961.6667
$ws.Range("K132").Value = 1476
$ws.Range("L132").Value = 2885.0001
$ws.Range("M132").Value = 1054
$ws.Range("N132").Value = -7945.0001
